$d = $word.ActiveDocument

# 1) Remove the existing "_GoBack" bookmark (located in an empty table cell paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2) Collapse the struck-through / highlighted date runs back to their plain originals.

# "Procure item(s) by <strike>May 31, 2020</strike> June 15, 2020 to support..."
$found = $d.Content.Find.Execute(
    "Procure item(s) by May 31, 2020 June 15, 2020 to support violence prevention programming.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Procure item(s) by May 31, 2020 to support violence prevention programming.",
    2)

# "If applicable, install item(s) by <strike>May 31, 2020</strike> June 15, 2020."
$found = $d.Content.Find.Execute(
    "If applicable, install item(s) by May 31, 2020 June 15, 2020.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "If applicable, install item(s) by May 31, 2020.",
    2)

# "Begin using item(s) in current violence prevention programming by <strike>June 1, 2020</strike> June 20, 2020."
$found = $d.Content.Find.Execute(
    "Begin using item(s) in current violence prevention programming by June 1, 2020 June 20, 2020.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Begin using item(s) in current violence prevention programming by June 1, 2020.",
    2)

# 3) Add a fresh "_GoBack" bookmark at the end of the "Describe how the one time-supports..." paragraph.
$target = $d.Content.Find.Execute(
    "Describe how the one time-supports will be used through June 30, 2021. Note: Selected applicants will have to certify that the purchased items will be utilized for violence prevention activities through June 30, 2021. ",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$searchRange = $d.Content
$searchRange.Find.Execute(
    "Describe how the one time-supports will be used through June 30, 2021. Note: Selected applicants will have to certify that the purchased items will be utilized for violence prevention activities through June 30, 2021. ",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$bmRange = $d.Range($searchRange.End, $searchRange.End)
$d.Bookmarks.Add("_GoBack", $bmRange)
